$wb = $excel.ActiveWorkbook

# zh-cn sheet: update handoff/handback datetime strings on row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-12 06:45:30"
$wsZhCn.Range("H2").Value = "2016-03-12 06:45:47"

# de-de sheet: update handoff/handback datetime strings on row 2
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-12 06:45:33"
$wsDeDe.Range("H2").Value = "2016-03-12 06:45:52"
